$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new row's "Result" text (contains embedded newlines/tabs exactly
# as produced by the Selenium/Chrome error trace). The here-string strips the
# final newline before the closing '@, so a blank line is included to
# preserve the trailing newline from the original text.
$resultText = @'
Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=128.0.6613.138); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF7BD199412+29090]
	(No symbol) [0x00007FF7BD10E239]
	(No symbol) [0x00007FF7BCFCB1DA]
	(No symbol) [0x00007FF7BD01EFE7]
	(No symbol) [0x00007FF7BD01F23C]
	(No symbol) [0x00007FF7BD0697C7]
	(No symbol) [0x00007FF7BD04672F]
	(No symbol) [0x00007FF7BD0665A2]
	(No symbol) [0x00007FF7BD046493]
	(No symbol) [0x00007FF7BD0109D1]
	(No symbol) [0x00007FF7BD011B31]
	GetHandleVerifier [0x00007FF7BD4B871D+3302573]
	GetHandleVerifier [0x00007FF7BD504243+3612627]
	GetHandleVerifier [0x00007FF7BD4FA417+3572135]
	GetHandleVerifier [0x00007FF7BD255EB6+801862]
	(No symbol) [0x00007FF7BD11945F]
	(No symbol) [0x00007FF7BD114FB4]
	(No symbol) [0x00007FF7BD115140]
	(No symbol) [0x00007FF7BD10461F]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]

'@

$row = 14
$ws.Cells.Item($row, 1).Value = "2024-09-23 17:03:50"
$ws.Cells.Item($row, 2).Value = "check_availability"
$ws.Cells.Item($row, 3).Value = "https://www.opentable.com/r/hals-the-steakhouse-nashville"
$ws.Cells.Item($row, 4).Value = $resultText

# Column E holds a date-looking string ("2024-09-23") that must stay plain
# text (matching the other rows) rather than being auto-converted into a
# date serial number by Excel. Force text interpretation, then drop back to
# the default "Normal" style so no explicit cell style is left behind.
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "2024-09-23"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).Value = "17:03:50"

# Setting the long, multi-line Result text auto-expands the row height;
# restore the default (non-custom) row height to match the original layout.
$ws.Rows.Item($row).EntireRow.AutoFit()
